$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 384
    3 = 160
    4 = 32
    5 = 347
    6 = 347
    7 = 85
    8 = 49
    9 = 195
    10 = 33
    11 = 102
    12 = 715
    13 = 57
    14 = 58
    15 = 506
    16 = 506
    17 = 317
    18 = 285
    19 = 34
    20 = 43
    21 = 103
    22 = 46
    23 = 100
    24 = 53
    25 = 47
    26 = 116
    27 = 65
    28 = 33
    29 = 26
    30 = 44
    31 = 44
    32 = 58
    33 = 67
    34 = 20
    35 = 49
    36 = 44
    37 = 25
    38 = 27
    39 = 160
    40 = 465
    41 = 441
    42 = 29
    43 = 501
    44 = 501
    45 = 51
    46 = 26
    47 = 56
    48 = 30
    49 = 345
    50 = 345
    51 = 34
    52 = 250
    53 = 35
    54 = 45
    55 = 82
    56 = 26
    57 = 336
    58 = 336
    59 = 264
    60 = 30
    61 = 719
    62 = 31
    63 = 30
    64 = 40
    65 = 59
    66 = 29
    67 = 49
    68 = 34
    69 = 31
    70 = 26
    71 = 35
    72 = 490
    73 = 343
    74 = 25
    75 = 24
    76 = 79
    77 = 56
    78 = 536
    79 = 536
    80 = 36
    81 = 25
    82 = 40
    83 = 67
    84 = 515
    85 = 46
    86 = 34
    87 = 32
    88 = 34
    89 = 30
    90 = 27
    91 = 30
    92 = 85
    93 = 72
    94 = 494
    95 = 69
    96 = 179
    97 = 53
    98 = 30
    99 = 41
    100 = 649
    101 = 649
    102 = 471
    103 = 51
    104 = 134
    105 = 36
    106 = 72
    107 = 32
    108 = 26
    109 = 32
    110 = 29
    111 = 27
    112 = 65
    113 = 83
    114 = 111
    115 = 734
    116 = 34
    117 = 45
    118 = 136
    119 = 47
    120 = 28
    121 = 26
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
